$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 0.09425133333333334
$ws.Range("H2").Value = 0.282754
$ws.Range("I2").Value = 0.02715992817009031
$ws.Range("J2").Value = 0.02715992817009031
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2807486666666667
$ws.Range("N2").Value = 0.8422460000000001
$ws.Range("O2").Value = 0.03211396410631209
$ws.Range("P2").Value = 0.03211396410631208
$ws.Range("Q2").Value = 0.0264609361648889
$ws.Range("R2").Value = 0.238148425484
$ws.Range("S2").Value = 0.0008722129583842947
$ws.Range("T2").Value = 0.0008722129583842944

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 0.09425133333333334
$ws.Range("H3").Value = 0.282754
$ws.Range("I3").Value = 0.02715992817009031
$ws.Range("J3").Value = 0.02715992817009031
$ws.Range("O3").Value = 0.06996648921957034
$ws.Range("P3").Value = 0.06996648921957033
$ws.Range("Q3").Value = 0.05765027322044444
$ws.Range("R3").Value = 0.518852458984
$ws.Range("S3").Value = 0.001900284821516928
$ws.Range("T3").Value = 0.001900284821516928

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 0.09425133333333334
$ws.Range("H4").Value = 0.282754
$ws.Range("I4").Value = 0.02715992817009031
$ws.Range("J4").Value = 0.02715992817009031
$ws.Range("O4").Value = 0.8979195466741177
$ws.Range("P4").Value = 0.8979195466741176
$ws.Range("Q4").Value = 0.7398585776297779
$ws.Range("R4").Value = 6.658727198668001
$ws.Range("S4").Value = 0.02438743039018909
$ws.Range("T4").Value = 0.02438743039018908

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.357039508851706
$ws.Range("J5").Value = 0.357039508851706
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2807486666666667
$ws.Range("N5").Value = 0.8422460000000001
$ws.Range("O5").Value = 0.03211396410631209
$ws.Range("P5").Value = 0.03211396410631208
$ws.Range("Q5").Value = 0.3478506862353333
$ws.Range("R5").Value = 3.130656176118
$ws.Range("S5").Value = 0.01146595397179898
$ws.Range("T5").Value = 0.01146595397179898

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.357039508851706
$ws.Range("J6").Value = 0.357039508851706
$ws.Range("O6").Value = 0.06996648921957034
$ws.Range("P6").Value = 0.06996648921957033
$ws.Range("S6").Value = 0.02498080094703358
$ws.Range("T6").Value = 0.02498080094703357

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.357039508851706
$ws.Range("J7").Value = 0.357039508851706
$ws.Range("O7").Value = 0.8979195466741177
$ws.Range("P7").Value = 0.8979195466741176
$ws.Range("S7").Value = 0.3205927539328735
$ws.Range("T7").Value = 0.3205927539328734

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.6158005629782037
$ws.Range("J8").Value = 0.6158005629782037
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2807486666666667
$ws.Range("N8").Value = 0.8422460000000001
$ws.Range("O8").Value = 0.03211396410631209
$ws.Range("P8").Value = 0.03211396410631208
$ws.Range("Q8").Value = 0.5999522268697779
$ws.Range("R8").Value = 5.399570041828
$ws.Range("S8").Value = 0.01977579717612881
$ws.Range("T8").Value = 0.01977579717612881

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.6158005629782037
$ws.Range("J9").Value = 0.6158005629782037
$ws.Range("O9").Value = 0.06996648921957034
$ws.Range("P9").Value = 0.06996648921957033
$ws.Range("S9").Value = 0.04308540345101984
$ws.Range("T9").Value = 0.04308540345101983

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.6158005629782037
$ws.Range("J10").Value = 0.6158005629782037
$ws.Range("O10").Value = 0.8979195466741177
$ws.Range("P10").Value = 0.8979195466741176
$ws.Range("S10").Value = 0.5529393623510551
$ws.Range("T10").Value = 0.5529393623510551
